# Update the "Förändrad" (column C) date serial value from 45172 to 45175
# for data rows 2 through 28 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value = 45175
    }
}
